# Adding 10 Manhattan plots lecture and recitation
#
# Rebuilds the Schedule / Schedule_date tables with the updated week
# list: week 8's topic is filled in, the Manhattan-plots week absorbs
# "making lots of plots at once", a new ggplot-extensions week is added,
# a Thanksgiving break week is added, and a final capstone week is added.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Schedule_date")

# ---------------------------------------------------------------------
# Schedule: A=Week, B=Module, C=Topic  (16 data rows, was 15)
# ---------------------------------------------------------------------
$ws1.Rows.Item(15).Insert()

$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "1: Principles"
$ws1.Range("C2").Value = "Principles of data visualization"
$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "1: Principles"
$ws1.Range("C3").Value = "Good and bad visualizations"
$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = "2: Coding fundamentals"
$ws1.Range("C4").Value = "R Markdown for reproducible research"
$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = "2: Coding fundamentals"
$ws1.Range("C5").Value = "ggplot 101"
$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = "2: Coding fundamentals"
$ws1.Range("C6").Value = "Themes, labels, facets (ggplot 102)"
$ws1.Range("A7").Value = 6
$ws1.Range("B7").Value = "3: Data exploration"
$ws1.Range("C7").Value = "Data distributions"
$ws1.Range("A8").Value = 7
$ws1.Range("B8").Value = "3: Data exploration"
$ws1.Range("C8").Value = "Correlations"
$ws1.Range("A9").Value = 8
$ws1.Range("B9").Value = "Open session, capstone prep"
$ws1.Range("C9").Value = "Open session, capstone prep"
$ws1.Range("A10").Value = 9
$ws1.Range("B10").Value = "3: Data exploration"
$ws1.Range("C10").Value = "Annotating statistics"
$ws1.Range("A11").Value = 10
$ws1.Range("B11").Value = "4: Putting it together"
$ws1.Range("C11").Value = "Principal components analysis"
$ws1.Range("A12").Value = 11
$ws1.Range("B12").Value = "4: Putting it together"
$ws1.Range("C12").Value = "Manhattan plots and making lots of plots at once"
$ws1.Range("A13").Value = 12
$ws1.Range("B13").Value = "4: Putting it together"
$ws1.Range("C13").Value = "Interactive plots"
$ws1.Range("A14").Value = 13
$ws1.Range("B14").Value = "4: Putting it together"
$ws1.Range("C14").Value = "ggplot extension packages and complexheatmap"
$ws1.Range("A15").Value = 14
$ws1.Range("B15").Value = "No class, Thanksgiving"
$ws1.Range("C15").Value = "Relaxing and eating"
$ws1.Range("A16").Value = 15
$ws1.Range("B16").Value = "4: Putting it together"
$ws1.Range("C16").Value = "Capstone assignment open session"
$ws1.Range("A17").Value = 16
$ws1.Range("B17").Value = "4: Putting it together"
$ws1.Range("C17").Value = "Capstone assignment open session"

# Column B needs to widen to fit the longer module/topic text
$ws1.Columns.Item(2).ColumnWidth = 23.92
$ws1.Range("F23").Select()

# ---------------------------------------------------------------------
# Schedule_date: A=Week, B=Date, C=Module, D=Topic  (16 data rows, was 14)
# ---------------------------------------------------------------------
$ws2.Rows.Item(10).Insert()
$ws2.Rows.Item(17).Insert()

$ws2.Range("A2").Value = 1
$ws2.Range("B2").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B2").Value = [datetime]"2022-08-23"
$ws2.Range("C2").Value = "1: Principles"
$ws2.Range("D2").Value = "Principles of data visualization"
$ws2.Range("A3").Value = 2
$ws2.Range("B3").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B3").Value = [datetime]"2022-08-30"
$ws2.Range("C3").Value = "1: Principles"
$ws2.Range("D3").Value = "Good and bad visualizations"
$ws2.Range("A4").Value = 3
$ws2.Range("B4").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B4").Value = [datetime]"2022-09-06"
$ws2.Range("C4").Value = "2: Coding fundamentals"
$ws2.Range("D4").Value = "R Markdown for reproducible research"
$ws2.Range("A5").Value = 4
$ws2.Range("B5").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B5").Value = [datetime]"2022-09-13"
$ws2.Range("C5").Value = "2: Coding fundamentals"
$ws2.Range("D5").Value = "ggplot 101"
$ws2.Range("A6").Value = 5
$ws2.Range("B6").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B6").Value = [datetime]"2022-09-21"
$ws2.Range("C6").Value = "2: Coding fundamentals"
$ws2.Range("D6").Value = "Themes, labels, facets (ggplot 102)"
$ws2.Range("A7").Value = 6
$ws2.Range("B7").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B7").Value = [datetime]"2022-09-28"
$ws2.Range("C7").Value = "3: Data exploration"
$ws2.Range("D7").Value = "Data distributions"
$ws2.Range("A8").Value = 7
$ws2.Range("B8").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B8").Value = [datetime]"2022-10-04"
$ws2.Range("C8").Value = "3: Data exploration"
$ws2.Range("D8").Value = "Correlations"
$ws2.Range("A9").Value = 8
$ws2.Range("B9").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B9").Value = [datetime]"2022-10-11"
$ws2.Range("C9").Value = "Open session, capstone prep"
$ws2.Range("D9").Value = "Open session, capstone prep"
$ws2.Range("A10").Value = 9
$ws2.Range("B10").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B10").Value = [datetime]"2022-10-18"
$ws2.Range("C10").Value = "3: Data exploration"
$ws2.Range("D10").Value = "Annotating statistics"
$ws2.Range("A11").Value = 10
$ws2.Range("B11").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B11").Value = [datetime]"2022-10-25"
$ws2.Range("C11").Value = "4: Putting it together"
$ws2.Range("D11").Value = "Principal components analysis"
$ws2.Range("A12").Value = 11
$ws2.Range("B12").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B12").Value = [datetime]"2022-11-01"
$ws2.Range("C12").Value = "4: Putting it together"
$ws2.Range("D12").Value = "Manhattan plots and making lots of plots at once"
$ws2.Range("A13").Value = 12
$ws2.Range("B13").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B13").Value = [datetime]"2022-11-08"
$ws2.Range("C13").Value = "4: Putting it together"
$ws2.Range("D13").Value = "Interactive plots"
$ws2.Range("A14").Value = 13
$ws2.Range("B14").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B14").Value = [datetime]"2022-11-15"
$ws2.Range("C14").Value = "4: Putting it together"
$ws2.Range("D14").Value = "ggplot extension packages and complexheatmap"
$ws2.Range("A15").Value = 14
$ws2.Range("B15").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B15").Value = [datetime]"2022-11-22"
$ws2.Range("C15").Value = "No class, Thanksgiving"
$ws2.Range("D15").Value = "Relaxing and eating"
$ws2.Range("A16").Value = 15
$ws2.Range("B16").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B16").Value = [datetime]"2022-11-29"
$ws2.Range("C16").Value = "4: Putting it together"
$ws2.Range("D16").Value = "Capstone assignment open session"
$ws2.Range("A17").Value = 16
$ws2.Range("B17").NumberFormat = "d\-mmm\-yy"
$ws2.Range("B17").Value = [datetime]"2022-12-05"
$ws2.Range("C17").Value = "4: Putting it together"
$ws2.Range("D17").Value = "Capstone assignment open session"

$ws2.Range("A1:D17").Select()

$ws1.Activate()
